$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now represents the Gold futures (GC=F) instrument; row 3 now
# represents the GLD ETF instrument - the two rows effectively swapped
# identity/content versus the previous snapshot, and several metrics
# were refreshed with new values.

# Row 2 (Gold Feb 26 / GC=F)
$ws.Range("B2").Value = "Gold Feb 26"
$ws.Range("C2").Value = "GC=F"
$ws.Range("D2").Value = 4237.1
$ws.Range("E2").Value = 55.6
$ws.Range("F2").Value = 1.73
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 83
$ws.Range("K2").Value = 67.7
$ws.Range("N2").Value = 54.86376272656823

# Row 3 (StreetTRACKS Gold Shares / GLD)
$ws.Range("B3").Value = "StreetTRACKS Gold Shares"
$ws.Range("C3").Value = "GLD"
$ws.Range("D3").Value = 387.13
$ws.Range("E3").Value = 56.3
$ws.Range("F3").Value = 1.05
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 73
$ws.Range("I3").Value = 83
$ws.Range("J3").Value = 96
$ws.Range("K3").Value = 67.7
$ws.Range("N3").Value = 54.86376272656823

# Row 4 (Newmont Corporation / NEM) - only the macro score refreshed
$ws.Range("N4").Value = 54.86376272656823
